# Auto-generated edits applying the betting-odds/correct-score updates for
# Jogos_da_Semana_FlashScore_2025-02-11.xlsx (rows 2-27 of Sheet1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("Q2").Value = 2.17
$ws.Range("R2").Value = 1.7

# Row 3
$ws.Range("J3").Value = 2.38
$ws.Range("Q3").Value = 1.95
$ws.Range("R3").Value = 1.9
$ws.Range("S3").Value = 2.5
$ws.Range("T3").Value = 1.5
$ws.Range("W3").Value = 5
$ws.Range("X3").Value = 1.17
$ws.Range("AN3").Value = 11
$ws.Range("AO3").Value = 29

# Row 4
$ws.Range("G4").Value = 1.62
$ws.Range("H4").Value = 3.7
$ws.Range("I4").Value = 6
$ws.Range("Q4").Value = 1.74
$ws.Range("R4").Value = 2.11
$ws.Range("U4").Value = 3.6
$ws.Range("V4").Value = 1.3
$ws.Range("AH4").Value = 34
$ws.Range("AJ4").Value = 7
$ws.Range("AK4").Value = 21
$ws.Range("AP4").Value = 19

# Row 5
$ws.Range("H5").Value = 3.5
$ws.Range("J5").Value = 2.5
$ws.Range("N5").Value = 7.5
$ws.Range("AI5").Value = 7.5
$ws.Range("AK5").Value = 21
$ws.Range("AN5").Value = 10

# Row 6
$ws.Range("G6").Value = 3.6
$ws.Range("H6").Value = 3
$ws.Range("K6").Value = 2
$ws.Range("O6").Value = 1.4
$ws.Range("P6").Value = 2.75
$ws.Range("Q6").Value = 1.69
$ws.Range("R6").Value = 2.19
$ws.Range("S6").Value = 2.35
$ws.Range("T6").Value = 1.57
$ws.Range("U6").Value = 3.4
$ws.Range("V6").Value = 1.32
$ws.Range("W6").Value = 4.33
$ws.Range("X6").Value = 1.2
$ws.Range("Y6").Value = 1.5
$ws.Range("Z6").Value = 2.5
$ws.Range("AA6").Value = 1.95
$ws.Range("AB6").Value = 1.8
$ws.Range("AC6").Value = 9
$ws.Range("AF6").Value = 41
$ws.Range("AG6").Value = 34
$ws.Range("AI6").Value = 7.5
$ws.Range("AJ6").Value = 5.5
$ws.Range("AM6").Value = 351
$ws.Range("AN6").Value = 6.5
$ws.Range("AR6").Value = 21
$ws.Range("AS6").Value = 34

# Row 7
$ws.Range("G7").Value = 1.23
$ws.Range("H7").Value = 5.9
$ws.Range("I7").Value = 9.25
$ws.Range("J7").Value = 1.6
$ws.Range("K7").Value = 2.77
$ws.Range("L7").Value = 7.3
$ws.Range("AA7").Value = 1.78
$ws.Range("AB7").Value = 1.82
$ws.Range("AD7").Value = 7.5
$ws.Range("AE7").Value = 9.25
$ws.Range("AF7").Value = 8
$ws.Range("AG7").Value = 10
$ws.Range("AH7").Value = 24
$ws.Range("AJ7").Value = 12.5
$ws.Range("AK7").Value = 22
$ws.Range("AL7").Value = 80
$ws.Range("AM7").Value = 500
$ws.Range("AN7").Value = 32
$ws.Range("AO7").Value = 75
$ws.Range("AP7").Value = 30
$ws.Range("AQ7").Value = 250
$ws.Range("AR7").Value = 100
$ws.Range("AS7").Value = 75

# Row 8
$ws.Range("Q8").Value = 1.63
$ws.Range("R8").Value = 2.24
$ws.Range("U8").Value = 3.2
$ws.Range("V8").Value = 1.34

# Row 10
$ws.Range("G10").Value = 2.8
$ws.Range("H10").Value = 2.92
$ws.Range("I10").Value = 2.62
$ws.Range("J10").Value = 3.55
$ws.Range("K10").Value = 1.87
$ws.Range("L10").Value = 3.4
$ws.Range("M10").Value = 1.13
$ws.Range("N10").Value = 5.1
$ws.Range("O10").Value = 1.6
$ws.Range("P10").Value = 2.18
$ws.Range("S10").Value = 2.75
$ws.Range("T10").Value = 1.39
$ws.Range("W10").Value = 5.1
$ws.Range("X10").Value = 1.13
$ws.Range("Y10").Value = 1.6
$ws.Range("Z10").Value = 2.18
$ws.Range("AA10").Value = 2.27
$ws.Range("AB10").Value = 1.57
$ws.Range("AC10").Value = 6.1
$ws.Range("AD10").Value = 12
$ws.Range("AE10").Value = 11.5
$ws.Range("AF10").Value = 32
$ws.Range("AG10").Value = 32
$ws.Range("AH10").Value = 60
$ws.Range("AI10").Value = 5.1
$ws.Range("AJ10").Value = 5.9
$ws.Range("AK10").Value = 21
$ws.Range("AL10").Value = 150
$ws.Range("AM10").Value = 101
$ws.Range("AN10").Value = 5.8
$ws.Range("AO10").Value = 11
$ws.Range("AP10").Value = 11.25
$ws.Range("AQ10").Value = 29
$ws.Range("AR10").Value = 30
$ws.Range("AS10").Value = 60

# Row 13
$ws.Range("G13").Value = 1.8
$ws.Range("I13").Value = 4.5
$ws.Range("S13").Value = 1.88
$ws.Range("T13").Value = 1.98
$ws.Range("AI13").Value = 11
$ws.Range("AM13").Value = 201

# Row 14
$ws.Range("K14").Value = 2
$ws.Range("AJ14").Value = 6
$ws.Range("AL14").Value = 51
$ws.Range("AN14").Value = 9.5
$ws.Range("AO14").Value = 19

# Row 15
$ws.Range("S15").Value = 2.1
$ws.Range("T15").Value = 1.73

# Row 17
$ws.Range("G17").Value = 1.75
$ws.Range("I17").Value = 4.75
$ws.Range("J17").Value = 2.4
$ws.Range("S17").Value = 2
$ws.Range("T17").Value = 1.8
$ws.Range("AD17").Value = 8
$ws.Range("AI17").Value = 9

# Row 19
$ws.Range("I19").Value = 3.4
$ws.Range("M19").Value = 1.07
$ws.Range("N19").Value = 9
$ws.Range("AI19").Value = 9

# Row 20
$ws.Range("G20").Value = 4.5
$ws.Range("I20").Value = 1.75
$ws.Range("AJ20").Value = 7
$ws.Range("AK20").Value = 17
$ws.Range("AQ20").Value = 13

# Row 24
$ws.Range("AI24").Value = 12

# Row 26
$ws.Range("G26").Value = 3.2
$ws.Range("I26").Value = 2.2
$ws.Range("J26").Value = 3.75
$ws.Range("L26").Value = 2.88
$ws.Range("S26").Value = 1.93
$ws.Range("T26").Value = 1.88
$ws.Range("AA26").Value = 1.73
$ws.Range("AB26").Value = 2
$ws.Range("AC26").Value = 11
$ws.Range("AD26").Value = 17
$ws.Range("AE26").Value = 12
$ws.Range("AF26").Value = 34
$ws.Range("AG26").Value = 26
$ws.Range("AH26").Value = 34
$ws.Range("AO26").Value = 11
$ws.Range("AP26").Value = 9
$ws.Range("AR26").Value = 17

# Row 27
$ws.Range("G27").Value = 7.6
$ws.Range("H27").Value = 5.2
$ws.Range("I27").Value = 1.32
$ws.Range("J27").Value = 6.4
$ws.Range("K27").Value = 2.67
$ws.Range("L27").Value = 1.72
$ws.Range("M27").Value = 1.03
$ws.Range("N27").Value = 9.75
$ws.Range("O27").Value = 1.14
$ws.Range("P27").Value = 4.9
$ws.Range("S27").Value = 1.45
$ws.Range("T27").Value = 2.57
$ws.Range("W27").Value = 2.1
$ws.Range("X27").Value = 1.65
$ws.Range("Y27").Value = 1.25
$ws.Range("Z27").Value = 3.6
$ws.Range("AA27").Value = 1.75
$ws.Range("AB27").Value = 1.98
$ws.Range("AC27").Value = 26
$ws.Range("AD27").Value = 55
$ws.Range("AE27").Value = 23
$ws.Range("AF27").Value = 175
$ws.Range("AG27").Value = 75
$ws.Range("AH27").Value = 60
$ws.Range("AI27").Value = 9.75
$ws.Range("AJ27").Value = 10.75
$ws.Range("AK27").Value = 18.5
$ws.Range("AL27").Value = 70
$ws.Range("AM27").Value = 450
$ws.Range("AN27").Value = 9.5
$ws.Range("AO27").Value = 7.6
$ws.Range("AP27").Value = 8.75
$ws.Range("AQ27").Value = 9
$ws.Range("AR27").Value = 10.25
$ws.Range("AS27").Value = 22
